$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '28.605.12', '  -3.11%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.850.28', '  -3.54%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.002', '  -1.14%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '335.91', '  +3.19%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.002', '  -0.96%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4652', '  -3.35%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3917', '  -3.31%  ')
    ,@(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07889', '  -3.98%  ')
    ,@(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.9833', '  -2.56%  ')
    ,@(11, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '22.19', '  -5.19%  ')
    ,@(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.934.93', '  +0.75%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.848', '  -3.37%  ')
    ,@(14, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.011', '  -3.10%  ')
    ,@(15, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06826', '  -0.64%  ')
    ,@(16, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.003', '  -1.03%  ')
    ,@(17, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '87.62', '  -4.18%  ')
    ,@(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001012', '  -2.59%  ')
    ,@(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.02', '  -2.87%  ')
    ,@(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  -0.93%  ')
    ,@(21, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '28.613.48', '  -3.08%  ')
    ,@(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.398', '  -4.80%  ')
    ,@(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.26', '  -5.11%  ')
    ,@(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.123', '  -3.21%  ')
    ,@(25, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.103.52', '  -1.94%  ')
    ,@(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '153.30', '  -1.74%  ')
    ,@(27, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '6.207', '  -4.99%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.40', '  -3.10%  ')
    ,@(29, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.022', '  -3.53%  ')
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '117.49', '  -2.56%  ')
    ,@(31, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.9750', '  -4.05%  ')
    ,@(32, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09440', '  -2.05%  ')
    ,@(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.375', '  -4.28%  ')
    ,@(34, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.504', '  -1.57%  ')
    ,@(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.351', '  -1.52%  ')
    ,@(36, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06129', '  -2.79%  ')
    ,@(37, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02194', '  -3.93%  ')
    ,@(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.162', '  -1.99%  ')
    ,@(39, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5695', '  -4.06%  ')
    ,@(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '7.580', '  -4.02%  ')
    ,@(41, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.12', '  -5.33%  ')
    ,@(42, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1789', '  -3.05%  ')
    ,@(43, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.388', '  -3.47%  ')
    ,@(44, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.252', '  +0.45%  ')
    ,@(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '11.83', '  -4.60%  ')
    ,@(46, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5382', '  -3.19%  ')
    ,@(47, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07143', '  -4.38%  ')
    ,@(48, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.907', '  -1.73%  ')
    ,@(49, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '113.22', '  -4.23%  ')
    ,@(50, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '43.49', '  +3.20%  ')
    ,@(51, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  -1.13%  ')
)

foreach ($row in $data) {
    $r = [int]$row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
